$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# AXA ROW rates update: refresh the start date and settlement currency
$ws.Range("D2").Value = "2025-06-01"
$ws.Range("E2").Value = "GBP"

# Reset the view back to the top-left corner / A2 selection
$ws.Range("A2").Select()
